$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Legg inn siste eksamen (2023 - Vår) i rad 14
$ws.Range("B14").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-23-v.pdf)"
$ws.Range("C14").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-23-v-fasit.pdf)"
$ws.Range("A14").Value = "2023 - Vår"

# Flytt markøren til neste ledige rad, slik som i den lagrede filen
$ws.Range("A15").Select()
